# England data now updates automatically from government website
# Two new date columns (CJ, CK) are appended after the existing CI column,
# for every region row (1-14). New cells inherit the same visual style as
# the preceding CI column by copying it (value+format) before overwriting
# the value with the real data point.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting of the last existing data column (CI) into the two
# new columns so the new cells share the identical cell style (border,
# font, number format, centered alignment) used throughout the table.
$ws.Range("CI1:CI14").Copy($ws.Range("CJ1"))
$ws.Range("CI1:CI14").Copy($ws.Range("CK1"))

# New data values per region (row) for columns CJ and CK.
$ws.Range("CJ1").Value = 1071
$ws.Range("CK1").Value = 1076

$ws.Range("CJ2").Value = 325
$ws.Range("CK2").Value = 326

$ws.Range("CJ3").Value = 261
$ws.Range("CK3").Value = 261

$ws.Range("CJ4").Value = 877
$ws.Range("CK4").Value = 877

$ws.Range("CJ5").Value = 944
$ws.Range("CK5").Value = 946

$ws.Range("CJ6").Value = 1283
$ws.Range("CK6").Value = 1284

$ws.Range("CJ7").Value = 3961
$ws.Range("CK7").Value = 3968

$ws.Range("CJ8").Value = 338
$ws.Range("CK8").Value = 338

$ws.Range("CJ9").Value = 2002
$ws.Range("CK9").Value = 2005

$ws.Range("CJ10").Value = 2757
$ws.Range("CK10").Value = 2760

$ws.Range("CJ11").Value = 8
$ws.Range("CK11").Value = 8

$ws.Range("CJ12").Value = 54
$ws.Range("CK12").Value = 54

$ws.Range("CJ13").Value = 1671
$ws.Range("CK13").Value = 1673

$ws.Range("CJ14").Value = 6
$ws.Range("CK14").Value = 6

# Update the view so the newly-added columns are visible / selected, same
# as what Excel records after a user scrolls to and selects the new range.
$ws.Range("A1:CK14").Select() | Out-Null
